$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped/changed rows 133 and 134 ---
$ws.Range("B133").Value = 7513577
$ws.Range("F133").Value = "Estoril"
$ws.Range("G133").Value = "SC Farense"
$ws.Range("H133").Value = 4
$ws.Range("K133").Value = 2.15
$ws.Range("L133").Value = 3.6
$ws.Range("M133").Value = 3.2
$ws.Range("N133").Value = 1.833
$ws.Range("O133").Value = 4
$ws.Range("P133").Value = 3.8
$ws.Range("Q133").Value = -0.5
$ws.Range("R133").Value = 1.875
$ws.Range("S133").Value = 1.975
$ws.Range("T133").Value = 2.75
$ws.Range("U133").Value = 1.875
$ws.Range("V133").Value = 1.975
$ws.Range("W133").Value = 0.833
$ws.Range("Z133").Value = 0.875
$ws.Range("AB133").Value = 0.875
$ws.Range("AC133").Value = -1
$ws.Range("B134").Value = 7515550
$ws.Range("F134").Value = "Gil Vicente"
$ws.Range("G134").Value = "Boavista"
$ws.Range("H134").Value = 1
$ws.Range("K134").Value = 2.3
$ws.Range("L134").Value = 3.5
$ws.Range("M134").Value = 2.9
$ws.Range("N134").Value = 2.3
$ws.Range("O134").Value = 3.3
$ws.Range("P134").Value = 3
$ws.Range("Q134").Value = -0.25
$ws.Range("R134").Value = 2.05
$ws.Range("S134").Value = 1.8
$ws.Range("T134").Value = 2.5
$ws.Range("U134").Value = 2.05
$ws.Range("V134").Value = 1.8
$ws.Range("W134").Value = 1.3
$ws.Range("Z134").Value = 1.05
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 0.8

# --- Update swapped/changed rows 220 and 221 ---
$ws.Range("B220").Value = 6875478
$ws.Range("F220").Value = "Boavista"
$ws.Range("G220").Value = "Moreirense"
$ws.Range("H220").Value = 1
$ws.Range("J220").Value = "H"
$ws.Range("K220").Value = 2.6
$ws.Range("L220").Value = 3.25
$ws.Range("M220").Value = 2.75
$ws.Range("N220").Value = 3.5
$ws.Range("O220").Value = 3
$ws.Range("P220").Value = 2.3
$ws.Range("Q220").Value = 0.25
$ws.Range("R220").Value = 1.825
$ws.Range("S220").Value = 2.025
$ws.Range("T220").Value = 2
$ws.Range("U220").Value = 1.975
$ws.Range("V220").Value = 1.875
$ws.Range("W220").Value = 2.5
$ws.Range("X220").Value = -1
$ws.Range("Z220").Value = 0.825
$ws.Range("AC220").Value = 0.875
$ws.Range("B221").Value = 6876651
$ws.Range("F221").Value = "Rio Ave"
$ws.Range("G221").Value = "Braga"
$ws.Range("H221").Value = 0
$ws.Range("J221").Value = "D"
$ws.Range("K221").Value = 4.333
$ws.Range("L221").Value = 4
$ws.Range("M221").Value = 1.727
$ws.Range("N221").Value = 4.5
$ws.Range("O221").Value = 4
$ws.Range("P221").Value = 1.7
$ws.Range("Q221").Value = 0.75
$ws.Range("R221").Value = 1.95
$ws.Range("S221").Value = 1.95
$ws.Range("T221").Value = 2.5
$ws.Range("U221").Value = 1.8
$ws.Range("V221").Value = 2.05
$ws.Range("W221").Value = -1
$ws.Range("X221").Value = 3
$ws.Range("Z221").Value = 0.95
$ws.Range("AC221").Value = 1.05

# --- Update rows 235 and 236 with new match data ---
$ws.Range("B235").Value = 6876668
$ws.Range("E235").Value = 45382.58333333334
$ws.Range("F235").Value = "Vizela"
$ws.Range("G235").Value = "Casa Pia"
$ws.Range("K235").Value = 2.375
$ws.Range("L235").Value = 3.1
$ws.Range("M235").Value = 3.1
$ws.Range("N235").Value = 2.375
$ws.Range("O235").Value = 3.1
$ws.Range("P235").Value = 3.3
$ws.Range("Q235").Value = -0.25
$ws.Range("R235").Value = 2.01
$ws.Range("S235").Value = 1.89
$ws.Range("T235").Value = 2.25
$ws.Range("U235").Value = 2.05
$ws.Range("V235").Value = 1.8
$ws.Range("B236").Value = 6876666
$ws.Range("E236").Value = 45383.67708333334
$ws.Range("F236").Value = "Portimonense"
$ws.Range("G236").Value = "Braga"
$ws.Range("K236").Value = 6
$ws.Range("L236").Value = 4.333
$ws.Range("M236").Value = 1.5
$ws.Range("N236").Value = 7.5
$ws.Range("O236").Value = 4.5
$ws.Range("P236").Value = 1.444
$ws.Range("Q236").Value = 1.25
$ws.Range("R236").Value = 1.89
$ws.Range("S236").Value = 2.01
$ws.Range("T236").Value = 3
$ws.Range("U236").Value = 2.025
$ws.Range("V236").Value = 1.825

# --- Remove trailing rows 237-241 (matches no longer present) ---
$ws.Range("A237:A241").EntireRow.Delete()
